$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in column C for rows 3..12 (alternating 10/1 -> 1/10)
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 10

# Update the view: scroll back to top (remove topLeftCell="A7") and change selection to C12
$ws.Activate()
$ws.Range("C12").Select()
